# Data refresh: the oldest event (row 2, "合肥·lovelive only", which already
# ended on 2024-07-20) was dropped from the scrape, so every following row
# shifts up by one. A handful of "想去人数" (interest-count) values were also
# refreshed to their latest scraped numbers.
#
# This affects the two worksheets that list every exhibition-type event:
#   "展览"     (sheet1) - exhibitions only
#   "全部类型" (sheet4) - all event types
# "演出" and "本地生活" do not contain that event and are left untouched.

$wb = $excel.ActiveWorkbook

# Updated "想去人数" (F column) values, keyed by the event's NEW row number
# (i.e. after row 2 has already been deleted and everything shifted up).
$sheet1Updates = @{
    2  = 8691
    4  = 37
    9  = 473
    10 = 80
    11 = 89
    13 = 6202
    14 = 209
    15 = 316
    16 = 2351
    17 = 112
    18 = 192
    20 = 466
}

$sheet4Updates = @{
    2  = 8691
    4  = 37
    11 = 473
    12 = 80
    13 = 89
    16 = 6202
    18 = 209
    19 = 316
    20 = 2351
    21 = 112
    22 = 192
    24 = 466
}

$targets = @(
    @{ Name = "展览";     Updates = $sheet1Updates; LastRow = 20 },
    @{ Name = "全部类型"; Updates = $sheet4Updates; LastRow = 24 }
)

foreach ($target in $targets) {
    $ws = $wb.Worksheets.Item($target.Name)

    # Remove the stale first event; everything below shifts up one row and
    # the sheet's used range shrinks by one row (dimension goes I..N+1 -> I..N).
    $ws.Rows("2:2").Delete()

    # Column A holds a literal running index (1, 2, 3, ...), not a formula,
    # so it needs to be renumbered after the shift.
    for ($r = 2; $r -le $target.LastRow; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }

    # Refresh the "想去人数" counts that changed since the last scrape.
    foreach ($row in $target.Updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $target.Updates[$row]
    }
}
